$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.865495324134827
$ws.Range("B1").Value = 3.207775115966797
$ws.Range("C1").Value = 2.522902250289917
$ws.Range("D1").Value = 2.411839962005615
$ws.Range("E1").Value = 2.295898914337158
